# Apply the workflow-slide edits described by the diff:
#  - shrink the "TextBox 7" placeholder height (cy 4524315 -> 4247317 EMU)
#  - fix "L&Ls" typo -> "L&L" in the binary-search bullet
#  - rewrite the "Find closest value..." bullet's text
#  - replace the "If it's within 0.5 mile..." bullet (promoted one indent
#    level) with the new "Create a new array..." bullet text
#  - delete the trailing "If not, end binary search function" bullet

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

# --- resize the shape (EMU -> points: 1 pt = 12700 EMU) ---
$shp.Height = 4247317 / 12700

$tr = $shp.TextFrame2.TextRange
$paragraphs = @($tr.Paragraphs())

# --- "Use binary search to find the closest L&Ls to the target value" ---
$runs = @($paragraphs[4].Runs())
$runs[0].Text = "Use binary search to find the closest L&L to the target value"

# --- "Find closest value, then compare to new address L&L" ---
$runs = @($paragraphs[5].Runs())
$runs[0].Text = "Once the closest value is found, iterate up and down in the JSON object"

# --- remove "If not, end binary search function" entirely ---
$paragraphs[7].Delete()
$paragraphs = @($tr.Paragraphs())

# --- remove "If it's within 0.5 mile of both L&L, add it to return object" ---
$paragraphs[6].Delete()
$paragraphs = @($tr.Paragraphs())

# --- insert the replacement bullet right after "Once the closest value..."
#     (paragraph index 5), splitting with a carriage return so the new
#     paragraph inherits that paragraph's indent level (lvl=1 / marL=800100)
[void]$paragraphs[5].InsertAfter("`rCreate a new array containing all of the values that are within the latitude equivalent of 0.5 miles")
